$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text cells, written in the exact order that reproduces the original
# --- author's shared-string table ordering (first-use order determines the
# --- <sst> index assigned to each unique string).
$ws.Range("A2").Value = "Check for E Records"
$ws.Range("B3").Value = "Check for Earliest E Record"
$ws.Range("C4").Value = "Fetch Details of Earlliest E Record"
$ws.Range("C5").Value = "Fetch Details of Earlliest E Record -1"
$ws.Range("B6").Value = "Check Container History [R] for Earlliest E Record"
$ws.Range("B10").Value = "Check If Trucker Code Changed"
$ws.Range("G11").Value = "PIERE"
$ws.Range("H10").Value = "GPON"
$ws.Range("G9").Value = "LAYTI"
$ws.Range("H9").Value = "KRTD"
$ws.Range("E9").Value = "BareMove"
$ws.Range("F9").Value = "OG"
$ws.Range("F11").Value = "BO- Bare out"
$ws.Range("F10").Value = "BI- Bare In"
$ws.Range("H8").Value = "Trucker CD"
$ws.Range("H13").Value = "Trucker Code Changed"
$ws.Range("H14").Value = "10 Mins added to 'E-1' record"
$ws.Range("G8").Value = "Location"
$ws.Range("E8").Value = "IG"
$ws.Range("F16").Value = "FO - Full out"
$ws.Range("H17").Value = "AMPF"
$ws.Range("H16").Value = "PMXN"
$ws.Range("E16").Value = "Bare Move"

# Remaining cells that repeat an already-used string (order among these is
# irrelevant - they just reuse the shared-string index established above).
$ws.Range("G18").Value = "PIERE"
$ws.Range("H11").Value = "GPON"
$ws.Range("G10").Value = "LAYTI"
$ws.Range("G16").Value = "LAYTI"
$ws.Range("G17").Value = "LAYTI"
$ws.Range("F18").Value = "BO- Bare out"
$ws.Range("F17").Value = "BI- Bare In"
$ws.Range("H18").Value = "AMPF"

# --- Numeric cells in column I (plain counts first)
$ws.Range("I9").Value = 28
$ws.Range("I10").Value = 31
$ws.Range("I11").Value = 21

# --- Bold styling for the "header" row of each mini table. Applied before
# --- any NumberFormat changes so the bold-only cell style is allocated
# --- first (matching the target cellXfs order).
$ws.Range("E8").Font.Bold = $true
$ws.Range("F10").Font.Bold = $true
$ws.Range("G10").Font.Bold = $true
$ws.Range("H10").Font.Bold = $true
$ws.Range("I10").Font.Bold = $true

# --- Time values formatted as h:mm in column I
$ws.Range("I16").Value = 0.48749999999999999
$ws.Range("I16").NumberFormat = "h:mm"

$ws.Range("I17").Value = 0.6166666666666667
$ws.Range("I17").NumberFormat = "h:mm"
$ws.Range("I17").Font.Bold = $true

$ws.Range("F17").Font.Bold = $true
$ws.Range("G17").Font.Bold = $true
$ws.Range("H17").Font.Bold = $true

$ws.Range("I18").Value = 0.60972222222222217
$ws.Range("I18").NumberFormat = "h:mm"

# --- Column widths (values chosen so the stored, quantized width in the
# --- saved file lands as close as possible to the target bestFit widths)
$ws.Columns("B").ColumnWidth = 41
$ws.Columns("F").ColumnWidth = 10.666666666666666
$ws.Columns("H").ColumnWidth = 9

# --- Page setup
$ws.PageSetup.Orientation = 1

# --- Final selection matches the saved workbook's active cell
$ws.Range("C4").Select() | Out-Null
